$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4572.241
$ws.Range("I15").Value = 4572.241
$ws.Range("K15").Value = 13716.723
$ws.Range("M15").Value = -13547.723

$ws.Range("H132").Value = 1305.4839
$ws.Range("I132").Value = 1076.2593
$ws.Range("K132").Value = 3228.7779
$ws.Range("M132").Value = -698.7779

$ws.Range("H137").Value = 1053.2858
$ws.Range("I137").Value = 849.6667
$ws.Range("K137").Value = 2549.0001
$ws.Range("M137").Value = 0.9998999999997977

$ws.Range("H141").Value = 3101.3667
$ws.Range("I141").Value = 2449.0435
$ws.Range("K141").Value = 7347.130500000001
$ws.Range("M141").Value = -2167.130500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 19749.25
$ws.Range("J24").Value = 19749.25
$ws.Range("L24").Value = 19749.25
$ws.Range("N24").Value = -20497.25

$ws.Range("H32").Value = 3090.2207
$ws.Range("I32").Value = 2166.283
$ws.Range("J32").Value = 6354.8
$ws.Range("K32").Value = 2166.283
$ws.Range("L32").Value = 6354.8
$ws.Range("M32").Value = -1879.283
$ws.Range("N32").Value = -6928.8

$ws.Range("H74").Value = 1174.68
$ws.Range("I74").Value = 1085.4117
$ws.Range("J74").Value = 1364.375
$ws.Range("K74").Value = 1085.4117
$ws.Range("L74").Value = 1364.375
$ws.Range("M74").Value = -211.4117000000001
$ws.Range("N74").Value = -3112.375

$ws.Range("H77").Value = 1174.68
$ws.Range("I77").Value = 1085.4117
$ws.Range("J77").Value = 1364.375
$ws.Range("K77").Value = 5427.058500000001
$ws.Range("L77").Value = 6821.875
$ws.Range("M77").Value = -1059.058500000001
$ws.Range("N77").Value = -15557.875

$ws.Range("H100").Value = 19749.25
$ws.Range("J100").Value = 19749.25
$ws.Range("L100").Value = 19749.25
$ws.Range("N100").Value = -21913.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 251791
$ws.Range("I86").Value = 1720.5
$ws.Range("J86").Value = 501861.5
$ws.Range("K86").Value = 1720.5
$ws.Range("L86").Value = 501861.5
$ws.Range("M86").Value = -597.5
$ws.Range("N86").Value = -504107.5

$ws.Range("H89").Value = 251791
$ws.Range("I89").Value = 1720.5
$ws.Range("J89").Value = 501861.5
$ws.Range("K89").Value = 8602.5
$ws.Range("L89").Value = 2509307.5
$ws.Range("M89").Value = -2986.5
$ws.Range("N89").Value = -2520539.5

$ws.Range("H105").Value = 2366.5833
$ws.Range("I105").Value = 2366.5833
$ws.Range("K105").Value = 2366.5833
$ws.Range("M105").Value = -619.5832999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2290.4614
$ws.Range("I31").Value = 1800
$ws.Range("J31").Value = 2710.8572
$ws.Range("K31").Value = 1800
$ws.Range("L31").Value = 2710.8572
$ws.Range("M31").Value = -1505
$ws.Range("N31").Value = -3300.8572

$ws.Range("H34").Value = 2290.4614
$ws.Range("I34").Value = 1800
$ws.Range("J34").Value = 2710.8572
$ws.Range("K34").Value = 1800
$ws.Range("L34").Value = 2710.8572
$ws.Range("M34").Value = -1598
$ws.Range("N34").Value = -3114.8572

$ws.Range("H99").Value = 2654.3333
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2981.5
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2981.5
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5977.5

$ws.Range("H126").Value = 2654.3333
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2981.5
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 8944.5
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13884.5

$ws.Range("H132").Value = 3119.0588
$ws.Range("I132").Value = 1722.4
$ws.Range("K132").Value = 5167.200000000001
$ws.Range("M132").Value = -2637.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 3999.6667
$ws.Range("J62").Value = 3999.6667
$ws.Range("L62").Value = 11999.0001
$ws.Range("N62").Value = -13371.0001

$ws.Range("H64").Value = 3035.3333
$ws.Range("I64").Value = 1712
$ws.Range("J64").Value = 3300
$ws.Range("K64").Value = 5136
$ws.Range("L64").Value = 9900
$ws.Range("M64").Value = -4866
$ws.Range("N64").Value = -10440

$ws.Range("H65").Value = 3999.6667
$ws.Range("J65").Value = 3999.6667
$ws.Range("L65").Value = 35997.0003
$ws.Range("N65").Value = -42861.0003

$ws.Range("H67").Value = 3035.3333
$ws.Range("I67").Value = 1712
$ws.Range("J67").Value = 3300
$ws.Range("K67").Value = 5136
$ws.Range("L67").Value = 9900
$ws.Range("M67").Value = -4200
$ws.Range("N67").Value = -11772

$ws.Range("H99").Value = 3200
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 9000
$ws.Range("L99").Value = 9750
$ws.Range("M99").Value = -6754
$ws.Range("N99").Value = -14242

$ws.Range("H106").Value = 4992
$ws.Range("I106").Value = 480
$ws.Range("K106").Value = 1440
$ws.Range("M106").Value = -494

$ws.Range("H131").Value = 775.3917
$ws.Range("J131").Value = 792.1539
$ws.Range("L131").Value = 2376.4617
$ws.Range("N131").Value = -12456.4617

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3186.0667
$ws.Range("I7").Value = 2480.2222
$ws.Range("J7").Value = 4244.8335
$ws.Range("K7").Value = 2480.2222
$ws.Range("L7").Value = 4244.8335
$ws.Range("M7").Value = -2368.2222
$ws.Range("N7").Value = -4468.8335

$ws.Range("H126").Value = 3186.0667
$ws.Range("I126").Value = 2480.2222
$ws.Range("J126").Value = 4244.8335
$ws.Range("K126").Value = 7440.6666
$ws.Range("L126").Value = 12734.5005
$ws.Range("M126").Value = -4970.6666
$ws.Range("N126").Value = -17674.5005

$ws.Range("H136").Value = 4340.75
$ws.Range("I136").Value = 3200
$ws.Range("J136").Value = 4911.125
$ws.Range("K136").Value = 9600
$ws.Range("L136").Value = 14733.375
$ws.Range("M136").Value = -7050
$ws.Range("N136").Value = -19833.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1527.1
$ws.Range("I132").Value = 1118.4546
$ws.Range("J132").Value = 2026.5555
$ws.Range("K132").Value = 3355.3638
$ws.Range("L132").Value = 6079.666499999999
$ws.Range("M132").Value = -825.3638000000001
$ws.Range("N132").Value = -11139.6665

$ws.Range("H136").Value = 79368264
$ws.Range("J136").Value = 3799.5
$ws.Range("L136").Value = 11398.5
$ws.Range("N136").Value = -16498.5
